{"js": "// Fix the misspelling \"comnpletely\" -> \"completely\" in the \"Don't\n// comnpletely commit to a tech stack\" sentence (Mitigations row of the\n// Risks table). The rendered sentence keeps reading \"Don't completely\n// commit to a tech stack\" - only the typo is corrected and the single\n// space between the two words is relocated accordingly.\n//\n// In the original OOXML the sentence tail is split across two runs:\n//   run1: \"comnpletely\"\n//   run2: \" commit to a tech stack\"   (leading space)\n// After the fix the split is:\n//   run1: \"completely \"               (trailing space)\n//   run2: \"commit to a tech stack\"\n// Office.js has no \"run\" object, so each of the two text spans is located\n// with body.search() (both needles are unique in the document) and\n// rewritten in place with Range.insertText(..., \"Replace\"), which keeps\n// each hit's run formatting (font, proofErr wrapping, etc.) intact.\n\nconst misspelled = context.document.body.search(\"comnpletely\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nmisspelled.load(\"items\");\nawait context.sync();\n\nif (misspelled.items.length === 0) {\n  throw new Error(\"Could not find misspelled word to fix\");\n}\nmisspelled.items[0].insertText(\"completely \", \"Replace\");\nawait context.sync();\n\nconst spaceBefore = context.document.body.search(\" commit to a tech stack\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nspaceBefore.load(\"items\");\nawait context.sync();\n\nif (spaceBefore.items.length === 0) {\n  throw new Error(\"Could not find ' commit to a tech stack' to fix\");\n}\nspaceBefore.items[0].insertText(\"commit to a tech stack\", \"Replace\");\nawait context.sync();\n", "ps1": "# Fix the misspelling \"comnpletely\" -> \"completely\" in the\n# \"Don't comnpletely commit to a tech stack\" sentence (Mitigations row of\n# the Risks table), and shift the space that originally sat in front of\n# \"commit\" to sit after \"completely\" instead, so the sentence still reads\n# \"Don't completely commit to a tech stack\".\n\n$d = $word.ActiveDocument\n\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Text = \"comnpletely\"\n$find1.Replacement.ClearFormatting()\n$find1.Replacement.Text = \"completely \"\n$find1.Execute([ref]$find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find1.Replacement.Text, 2)\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \" commit to a tech stack\"\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"commit to a tech stack\"\n$find2.Execute([ref]$find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find2.Replacement.Text, 2)\n"}
